# "re-added 1D noise results" -- rebuild the 1d_supp results table:
#   - drop the old blank spacer column (between "# layers" and "Activation Layer")
#   - relabel the remaining columns as Learning Rate / Momentum / Dropout (%)
#   - replace the old per-layer/activation values with the learning-rate &
#     momentum values for each uncertainty method
#   - rename "Direct" -> "Direct Regression"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1d_supp")

# The sheet was A:E with an empty column C acting as a spacer. Deleting it
# shifts the old D/E (Activation Layer / Dropout (%)) columns left into C/D,
# leaving a clean A:D range.
$ws.Columns.Item(3).Delete()

# Header row: B = Learning Rate, C = Momentum, D stays Dropout (%)
$ws.Cells.Item(1, 2).Value = "Learning Rate"
$ws.Cells.Item(1, 3).Value = "Momentum"

# Row labels
$ws.Cells.Item(2, 1).Value = "Dropout"
$ws.Cells.Item(3, 1).Value = "Direct Regression"
$ws.Cells.Item(4, 1).Value = "Bagging"
$ws.Cells.Item(5, 1).Value = "HydraNet (no direct uncertainty)"
$ws.Cells.Item(6, 1).Value = "HydraNet"

# Learning Rate / Momentum values per row (Dropout (%) column (D) is left
# as-is: 3 for Dropout, "---" for the rest)
$ws.Cells.Item(2, 2).Value = 0.05
$ws.Cells.Item(2, 3).Value = 0.5

$ws.Cells.Item(3, 2).Value = 0.0001
$ws.Cells.Item(3, 3).Value = 0

$ws.Cells.Item(4, 2).Value = 0.01
$ws.Cells.Item(4, 3).Value = 0.9

$ws.Cells.Item(5, 2).Value = 0.01
$ws.Cells.Item(5, 3).Value = 0.9

$ws.Cells.Item(6, 2).Value = 0.01
$ws.Cells.Item(6, 3).Value = 0.1

# Widen the Learning Rate column to fit the new header/values, and restore
# the author's last selection on this sheet.
$ws.Columns.Item(2).ColumnWidth = 14.85
$ws.Range("C26").Select() | Out-Null
